$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.282.54"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.638.81"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'601.60"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").Value = "'145.86"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.585"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "2.638.66"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").Value = "'27.20"
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("D15").Value = "3.110.47"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "63.183.71"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "'0.0000144"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").Value = "2.619.88"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "'11.40"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").Value = "'340.98"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").Value = "'6.88"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'5.57"
$ws.Range("E24").Value = "  -3.34%  "
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("E26").Value = "  -3.33%  "
$ws.Range("D27").Value = "'8.71"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").Value = "'551.20"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("B29").Value = "SuiNetwork"
$ws.Range("C29").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D29").Value = "'1.51"
$ws.Range("E29").Value = "  -5.33%  "
$ws.Range("D30").Value = "'0.163"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").Value = "'7.82"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("D35").Value = "0.0₃0803"
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("D36").Value = "'5.22"
$ws.Range("E36").Value = "  +5.84%  "
$ws.Range("D37").Value = "'166.21"
$ws.Range("E37").Value = "  -5.29%  "
$ws.Range("D39").Value = "'0.405"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("D40").Value = "'18.98"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("E41").Value = "  +4.34%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "'167.90"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'22.44"
$ws.Range("E44").Value = "  +1.85%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "'3.73"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("D46").Value = "'0.0570"
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "'18.64"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").Value = "'1.75"
$ws.Range("E51").Value = "  +1.10%  "
